$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A77").Value = "2025-04-29 11:36:39"
$ws.Range("B77").Value = 234
